$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row for columns D (4), J (10), K (11), L (12), M (13), P (16)
# after re-sorting the weekly price records (Jengibre, Vega Monumental Concepcion).
$rowData = @{
    2 = @(44664, 50, 11000, 12000, 11600, 892)
    3 = @(44838, 40, 14000, 15000, 14500, 1115)
    4 = @(44810, 50, 11000, 12000, 11600, 892)
    5 = @(44377, 40, 14000, 15000, 14500, 1115)
    6 = @(44320, 50, 26000, 28000, 26800, 2062)
    7 = @(44523, 40, 15000, 16000, 15500, 1192)
    8 = @(44719, 50, 13000, 14000, 13400, 1031)
    9 = @(44462, 60, 14000, 15000, 14500, 1115)
    10 = @(44755, 40, 14000, 15000, 14500, 1115)
    11 = @(44467, 100, 13000, 14000, 13500, 1038)
    12 = @(44510, 40, 15000, 16000, 15500, 1192)
    13 = @(44691, 100, 12000, 13000, 12500, 962)
    14 = @(44383, 50, 15000, 16000, 15400, 1185)
    15 = @(44433, 100, 13000, 14000, 13500, 1038)
    16 = @(44435, 100, 13000, 14000, 13500, 1038)
    17 = @(44313, 50, 25000, 26000, 25600, 1969)
    18 = @(44159, 60, 30000, 32000, 31000, 2385)
    19 = @(44316, 50, 27000, 28000, 27400, 2108)
    20 = @(44334, 50, 26000, 28000, 27200, 2092)
    21 = @(44308, 50, 26000, 27000, 26400, 2031)
    22 = @(44355, 60, 18000, 20000, 19000, 1462)
    23 = @(44708, 50, 13000, 14000, 13600, 1046)
    24 = @(44488, 40, 16000, 17000, 16500, 1269)
    25 = @(44610, 50, 17000, 18000, 17400, 1338)
    26 = @(44327, 50, 24000, 25000, 24400, 1877)
    27 = @(44362, 40, 15000, 16000, 15500, 1192)
    28 = @(44503, 35, 15000, 16000, 15429, 1187)
    29 = @(44509, 100, 15000, 16000, 15500, 1192)
    30 = @(44775, 20, 12000, 13000, 12500, 962)
    31 = @(44782, 40, 13000, 14000, 13500, 1038)
    32 = @(44769, 50, 14000, 15000, 14600, 1123)
    33 = @(44705, 50, 10000, 11000, 10400, 800)
    34 = @(44777, 25, 13000, 14000, 13600, 1046)
    35 = @(44474, 40, 13000, 14000, 13500, 1038)
    36 = @(44761, 25, 14000, 15000, 14400, 1108)
    37 = @(44425, 60, 14000, 15000, 14500, 1115)
    38 = @(44264, 40, 30000, 32000, 31000, 2385)
    39 = @(44453, 50, 14000, 15000, 14600, 1123)
    40 = @(44883, 60, 14000, 15000, 14500, 1115)
    41 = @(44813, 50, 13000, 14000, 13400, 1031)
    42 = @(44819, 50, 13000, 14000, 13400, 1031)
    43 = @(44741, 50, 14000, 15000, 14400, 1108)
    44 = @(44350, 40, 23000, 25000, 24000, 1846)
}

foreach ($row in $rowData.Keys) {
    $vals = $rowData[$row]
    $ws.Cells.Item($row, 4).Value  = $vals[0]   # D - Fecha
    $ws.Cells.Item($row, 10).Value = $vals[1]   # J - Volumen
    $ws.Cells.Item($row, 11).Value = $vals[2]   # K - Precio minimo
    $ws.Cells.Item($row, 12).Value = $vals[3]   # L - Precio maximo
    $ws.Cells.Item($row, 13).Value = $vals[4]   # M - Precio promedio ponderado
    $ws.Cells.Item($row, 16).Value = $vals[5]   # P - Precio $/Kg
}
